$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column AG: previously numeric 1, now the text label "1a"
$ws.Range("AG2").Value = "1a"

# Row 2, column AH: previously blank, now a date (same day as AH3: 2021-10-10)
$ws.Range("AH2").Value = Get-Date -Year 2021 -Month 10 -Day 10 -Hour 0 -Minute 0 -Second 0

# Reflect the new scroll position / active selection saved in the workbook view
$ws.Activate()
$ws.Range("AG3").Select()
$excel.ActiveWindow.ScrollColumn = 18
